# Weekly refresh of the Jengibre (ginger) price series: a new weekly
# observation is inserted as row 27 (pushing the existing history down by
# one row), and the sheet's used range grows from A1:R128 to A1:R129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 27; this shifts rows
# 27..128 down to 28..129 (values, not just formatting, move with them).
$ws.Rows.Item(27).Insert()

# Populate the newly-inserted row 27 with the latest weekly record.
$ws.Cells.Item(27, 1).Value  = 9
$ws.Cells.Item(27, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(27, 3).Value  = "Metropolitana"
$ws.Cells.Item(27, 4).Value  = 45054
$ws.Cells.Item(27, 5).Value  = 13
$ws.Cells.Item(27, 6).Value  = 100114007
$ws.Cells.Item(27, 7).Value  = "Jengibre"
$ws.Cells.Item(27, 8).Value  = "Sin especificar"
$ws.Cells.Item(27, 9).Value  = "Primera"
$ws.Cells.Item(27, 10).Value = 430
$ws.Cells.Item(27, 11).Value = 17000
$ws.Cells.Item(27, 12).Value = 18000
$ws.Cells.Item(27, 13).Value = 17500
$ws.Cells.Item(27, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 1346
$ws.Cells.Item(27, 17).Value = 13
$ws.Cells.Item(27, 18).Value = "Hortaliza"
